$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 74
$ws.Cells.Item(74, 8).Value = 50518.13  # H74: 48454.875 -> 50518.13
$ws.Cells.Item(74, 9).Value = 56943.05  # I74: 54145.9 -> 56943.05
$ws.Cells.Item(74, 11).Value = 56943.05  # K74: 54145.9 -> 56943.05
$ws.Cells.Item(74, 13).Value = -56007.05  # M74: -53209.9 -> -56007.05

# Row 77
$ws.Cells.Item(77, 8).Value = 50518.13  # H77: 48454.875 -> 50518.13
$ws.Cells.Item(77, 9).Value = 56943.05  # I77: 54145.9 -> 56943.05
$ws.Cells.Item(77, 11).Value = 284715.25  # K77: 270729.5 -> 284715.25
$ws.Cells.Item(77, 13).Value = -280035.25  # M77: -266049.5 -> -280035.25

# Row 97
$ws.Cells.Item(97, 8).Value = 20405.334  # H97: 22264.455 -> 20405.334
$ws.Cells.Item(97, 9).Value = 3933.3333  # I97: 4300 -> 3933.3333
$ws.Cells.Item(97, 10).Value = 25896  # J97: 29001.125 -> 25896
$ws.Cells.Item(97, 11).Value = 11799.9999  # K97: 12900 -> 11799.9999
$ws.Cells.Item(97, 12).Value = 77688  # L97: 87003.375 -> 77688
$ws.Cells.Item(97, 13).Value = -11303.9999  # M97: -12404 -> -11303.9999
$ws.Cells.Item(97, 14).Value = -78680  # N97: -87995.375 -> -78680

# Row 98
$ws.Cells.Item(98, 8).Value = 10000  # H98: 5260.9165 -> 10000
$ws.Cells.Item(98, 9).Value = 0  # I98: 4713.2 -> 0
$ws.Cells.Item(98, 10).Value = 10000  # J98: 7999.5 -> 10000
$ws.Cells.Item(98, 11).Value = 0  # K98: 4713.2 -> 0
$ws.Cells.Item(98, 12).Value = 10000  # L98: 7999.5 -> 10000
$ws.Cells.Item(98, 13).ClearContents()  # M98: -3215.2 -> (removed)
$ws.Cells.Item(98, 14).Value = -12996  # N98: -10995.5 -> -12996

# Row 101
$ws.Cells.Item(101, 8).Value = 204.2  # H101: 196 -> 204.2
$ws.Cells.Item(101, 9).Value = 205.5  # I101: 195.4 -> 205.5
$ws.Cells.Item(101, 11).Value = 616.5  # K101: 586.2 -> 616.5
$ws.Cells.Item(101, 13).Value = 1005.5  # M101: 1035.8 -> 1005.5

# Row 111
$ws.Cells.Item(111, 8).Value = 16015.429  # H111: 14061 -> 16015.429
$ws.Cells.Item(111, 9).Value = 977  # I111: 857.6 -> 977
$ws.Cells.Item(111, 11).Value = 2931  # K111: 2572.8 -> 2931
$ws.Cells.Item(111, 13).Value = 136  # M111: 494.1999999999998 -> 136

# Row 116
$ws.Cells.Item(116, 8).Value = 6773.5  # H116: 7511.143 -> 6773.5
$ws.Cells.Item(116, 10).Value = 8157  # J116: 12122.25 -> 8157
$ws.Cells.Item(116, 12).Value = 8157  # L116: 12122.25 -> 8157
$ws.Cells.Item(116, 14).Value = -15041  # N116: -19006.25 -> -15041

# Row 122
$ws.Cells.Item(122, 8).Value = 10000  # H122: 5260.9165 -> 10000
$ws.Cells.Item(122, 9).Value = 0  # I122: 4713.2 -> 0
$ws.Cells.Item(122, 10).Value = 10000  # J122: 7999.5 -> 10000
$ws.Cells.Item(122, 11).Value = 0  # K122: 14139.6 -> 0
$ws.Cells.Item(122, 12).Value = 30000  # L122: 23998.5 -> 30000
$ws.Cells.Item(122, 13).ClearContents()  # M122: -11689.6 -> (removed)
$ws.Cells.Item(122, 14).Value = -34900  # N122: -28898.5 -> -34900

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Cells.Item(45, 8).Value = 4978.0356  # H45: 4851.207 -> 4978.0356
$ws.Cells.Item(45, 9).Value = 6599.421  # I45: 6334.45 -> 6599.421
$ws.Cells.Item(45, 11).Value = 6599.421  # K45: 6334.45 -> 6599.421
$ws.Cells.Item(45, 13).Value = -6222.421  # M45: -5957.45 -> -6222.421

# Row 74
$ws.Cells.Item(74, 8).Value = 5079.606  # H74: 8104.3447 -> 5079.606
$ws.Cells.Item(74, 9).Value = 1264.1852  # I74: 1333.625 -> 1264.1852
$ws.Cells.Item(74, 10).Value = 22249  # J74: 40603.8 -> 22249
$ws.Cells.Item(74, 11).Value = 1264.1852  # K74: 1333.625 -> 1264.1852
$ws.Cells.Item(74, 12).Value = 22249  # L74: 40603.8 -> 22249
$ws.Cells.Item(74, 13).Value = -390.1851999999999  # M74: -459.625 -> -390.1851999999999
$ws.Cells.Item(74, 14).Value = -23997  # N74: -42351.8 -> -23997

# Row 77
$ws.Cells.Item(77, 8).Value = 5079.606  # H77: 8104.3447 -> 5079.606
$ws.Cells.Item(77, 9).Value = 1264.1852  # I77: 1333.625 -> 1264.1852
$ws.Cells.Item(77, 10).Value = 22249  # J77: 40603.8 -> 22249
$ws.Cells.Item(77, 11).Value = 6320.925999999999  # K77: 6668.125 -> 6320.925999999999
$ws.Cells.Item(77, 12).Value = 111245  # L77: 203019 -> 111245
$ws.Cells.Item(77, 13).Value = -1952.925999999999  # M77: -2300.125 -> -1952.925999999999
$ws.Cells.Item(77, 14).Value = -119981  # N77: -211755 -> -119981

# Row 132
$ws.Cells.Item(132, 8).Value = 3294.5  # H132: 3547.861 -> 3294.5
$ws.Cells.Item(132, 9).Value = 3072.7144  # I132: 3341.3447 -> 3072.7144
$ws.Cells.Item(132, 11).Value = 9218.143199999999  # K132: 10024.0341 -> 9218.143199999999
$ws.Cells.Item(132, 13).Value = -6688.143199999999  # M132: -7494.034100000001 -> -6688.143199999999

# Row 133
$ws.Cells.Item(133, 8).Value = 121054  # H133: 139995.8 -> 121054
$ws.Cells.Item(133, 10).Value = 121054  # J133: 139995.8 -> 121054
$ws.Cells.Item(133, 12).Value = 121054  # L133: 139995.8 -> 121054
$ws.Cells.Item(133, 14).Value = -126114  # N133: -145055.8 -> -126114

$ws = $wb.Worksheets.Item("BSM")
# Row 35
$ws.Cells.Item(35, 8).Value = 20000  # H35: 19000 -> 20000
$ws.Cells.Item(35, 9).Value = 20000  # I35: 19000 -> 20000
$ws.Cells.Item(35, 11).Value = 20000  # K35: 19000 -> 20000
$ws.Cells.Item(35, 13).Value = -19690  # M35: -18690 -> -19690

# Row 99
$ws.Cells.Item(99, 8).Value = 90573.95  # H99: 94676.82000000001 -> 90573.95
$ws.Cells.Item(99, 10).Value = 150215.72  # J99: 175199.83 -> 150215.72
$ws.Cells.Item(99, 12).Value = 150215.72  # L99: 175199.83 -> 150215.72
$ws.Cells.Item(99, 14).Value = -153211.72  # N99: -178195.83 -> -153211.72

# Row 107
$ws.Cells.Item(107, 8).Value = 1716.2  # H107: 1648.1482 -> 1716.2
$ws.Cells.Item(107, 9).Value = 1505.45  # I107: 1338.7778 -> 1505.45
$ws.Cells.Item(107, 10).Value = 2137.7  # J107: 2266.889 -> 2137.7
$ws.Cells.Item(107, 11).Value = 1505.45  # K107: 1338.7778 -> 1505.45
$ws.Cells.Item(107, 12).Value = 2137.7  # L107: 2266.889 -> 2137.7
$ws.Cells.Item(107, 13).Value = 414.55  # M107: 581.2221999999999 -> 414.55
$ws.Cells.Item(107, 14).Value = -5977.7  # N107: -6106.889 -> -5977.7

# Row 134
$ws.Cells.Item(134, 8).Value = 1526.3055  # H134: 1534.25 -> 1526.3055
$ws.Cells.Item(134, 9).Value = 1514.6177  # I134: 1523.0294 -> 1514.6177
$ws.Cells.Item(134, 11).Value = 4543.8531  # K134: 4569.0882 -> 4543.8531
$ws.Cells.Item(134, 13).Value = -2008.8531  # M134: -2034.0882 -> -2008.8531

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Cells.Item(31, 8).Value = 42148.848  # H31: 35639.71 -> 42148.848
$ws.Cells.Item(31, 9).Value = 65960.81  # I31: 50682.57 -> 65960.81
$ws.Cells.Item(31, 11).Value = 65960.81  # K31: 50682.57 -> 65960.81
$ws.Cells.Item(31, 13).Value = -65665.81  # M31: -50387.57 -> -65665.81

# Row 33
$ws.Cells.Item(33, 8).Value = 30020.908  # H33: 27486.363 -> 30020.908
$ws.Cells.Item(33, 9).Value = 5070  # I33: 4477 -> 5070
$ws.Cells.Item(33, 10).Value = 44278.57  # J33: 46660.832 -> 44278.57
$ws.Cells.Item(33, 11).Value = 5070  # K33: 4477 -> 5070
$ws.Cells.Item(33, 12).Value = 44278.57  # L33: 46660.832 -> 44278.57
$ws.Cells.Item(33, 13).Value = -4691  # M33: -4098 -> -4691
$ws.Cells.Item(33, 14).Value = -45036.57  # N33: -47418.832 -> -45036.57

# Row 34
$ws.Cells.Item(34, 8).Value = 42148.848  # H34: 35639.71 -> 42148.848
$ws.Cells.Item(34, 9).Value = 65960.81  # I34: 50682.57 -> 65960.81
$ws.Cells.Item(34, 11).Value = 65960.81  # K34: 50682.57 -> 65960.81
$ws.Cells.Item(34, 13).Value = -65758.81  # M34: -50480.57 -> -65758.81

# Row 62
$ws.Cells.Item(62, 8).Value = 18666.666  # H62: 14622.5 -> 18666.666
$ws.Cells.Item(62, 9).Value = 1000  # I62: 1745 -> 1000
$ws.Cells.Item(62, 11).Value = 1000  # K62: 1745 -> 1000
$ws.Cells.Item(62, 13).Value = -376  # M62: -1121 -> -376

# Row 65
$ws.Cells.Item(65, 8).Value = 18666.666  # H65: 14622.5 -> 18666.666
$ws.Cells.Item(65, 9).Value = 1000  # I65: 1745 -> 1000
$ws.Cells.Item(65, 11).Value = 5000  # K65: 8725 -> 5000
$ws.Cells.Item(65, 13).Value = -1880  # M65: -5605 -> -1880

# Row 132
$ws.Cells.Item(132, 8).Value = 2129.9736  # H132: 2073.1765 -> 2129.9736
$ws.Cells.Item(132, 9).Value = 1943.919  # I132: 2073.1765 -> 1943.919
$ws.Cells.Item(132, 10).Value = 9014  # J132: 0 -> 9014
$ws.Cells.Item(132, 11).Value = 5831.757000000001  # K132: 6219.529500000001 -> 5831.757000000001
$ws.Cells.Item(132, 12).Value = 27042  # L132: 0 -> 27042
$ws.Cells.Item(132, 13).Value = -3301.757000000001  # M132: -3689.529500000001 -> -3301.757000000001
$ws.Cells.Item(132, 14).Value = -32102  # N132: None -> -32102

# Row 134
$ws.Cells.Item(134, 8).Value = 10767.75  # H134: 9812.807000000001 -> 10767.75
$ws.Cells.Item(134, 9).Value = 6484.7144  # I134: 5786.625 -> 6484.7144
$ws.Cells.Item(134, 11).Value = 19454.1432  # K134: 17359.875 -> 19454.1432
$ws.Cells.Item(134, 13).Value = -16919.1432  # M134: -14824.875 -> -16919.1432

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Cells.Item(5, 8).Value = 758.6111  # H5: 733.2143 -> 758.6111
$ws.Cells.Item(5, 9).Value = 758.6111  # I5: 733.2143 -> 758.6111
$ws.Cells.Item(5, 11).Value = 2275.8333  # K5: 2199.6429 -> 2275.8333
$ws.Cells.Item(5, 13).Value = -2163.8333  # M5: -2087.6429 -> -2163.8333

# Row 131
$ws.Cells.Item(131, 8).Value = 26399.047  # H131: 25816.113 -> 26399.047
$ws.Cells.Item(131, 9).Value = 251611.75  # I131: 201439.4 -> 251611.75
$ws.Cells.Item(131, 11).Value = 754835.25  # K131: 604318.2 -> 754835.25
$ws.Cells.Item(131, 13).Value = -749795.25  # M131: -599278.2 -> -749795.25

# Row 135
$ws.Cells.Item(135, 8).Value = 758.6111  # H135: 733.2143 -> 758.6111
$ws.Cells.Item(135, 9).Value = 758.6111  # I135: 733.2143 -> 758.6111
$ws.Cells.Item(135, 11).Value = 6827.4999  # K135: 6598.928699999999 -> 6827.4999
$ws.Cells.Item(135, 13).Value = -4292.4999  # M135: -4063.928699999999 -> -4292.4999

# Row 137
$ws.Cells.Item(137, 8).Value = 2972.75  # H137: 2686.9285 -> 2972.75
$ws.Cells.Item(137, 9).Value = 1740  # I137: 1548 -> 1740
$ws.Cells.Item(137, 11).Value = 5220  # K137: 4644 -> 5220
$ws.Cells.Item(137, 13).Value = -120  # M137: 456 -> -120

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Cells.Item(70, 8).Value = 22534.908  # H70: 21073.75 -> 22534.908
$ws.Cells.Item(70, 9).Value = 12977  # I70: 11647.5 -> 12977
$ws.Cells.Item(70, 10).Value = 30499.834  # J70: 30500 -> 30499.834
$ws.Cells.Item(70, 11).Value = 12977  # K70: 11647.5 -> 12977
$ws.Cells.Item(70, 12).Value = 30499.834  # L70: 30500 -> 30499.834
$ws.Cells.Item(70, 13).Value = -12707  # M70: -11377.5 -> -12707
$ws.Cells.Item(70, 14).Value = -31039.834  # N70: -31040 -> -31039.834

# Row 73
$ws.Cells.Item(73, 8).Value = 22534.908  # H73: 21073.75 -> 22534.908
$ws.Cells.Item(73, 9).Value = 12977  # I73: 11647.5 -> 12977
$ws.Cells.Item(73, 10).Value = 30499.834  # J73: 30500 -> 30499.834
$ws.Cells.Item(73, 11).Value = 12977  # K73: 11647.5 -> 12977
$ws.Cells.Item(73, 12).Value = 30499.834  # L73: 30500 -> 30499.834
$ws.Cells.Item(73, 13).Value = -12041  # M73: -10711.5 -> -12041
$ws.Cells.Item(73, 14).Value = -32371.834  # N73: -32372 -> -32371.834

# Row 102
$ws.Cells.Item(102, 8).Value = 83334280  # H102: 90910110 -> 83334280
$ws.Cells.Item(102, 9).Value = 997.8889  # I102: 1052 -> 997.8889
$ws.Cells.Item(102, 10).Value = 333334140  # J102: 250000980 -> 333334140
$ws.Cells.Item(102, 11).Value = 997.8889  # K102: 1052 -> 997.8889
$ws.Cells.Item(102, 12).Value = 333334140  # L102: 250000980 -> 333334140
$ws.Cells.Item(102, 13).Value = 624.1111  # M102: 570 -> 624.1111
$ws.Cells.Item(102, 14).Value = -333337384  # N102: -250004224 -> -333337384

# Row 113
$ws.Cells.Item(113, 8).Value = 3266.4736  # H113: 3507.9412 -> 3266.4736
$ws.Cells.Item(113, 9).Value = 3358.3845  # I113: 3748.2727 -> 3358.3845
$ws.Cells.Item(113, 11).Value = 3358.3845  # K113: 3748.2727 -> 3358.3845
$ws.Cells.Item(113, 13).Value = -1188.3845  # M113: -1578.2727 -> -1188.3845

# Row 122
$ws.Cells.Item(122, 8).Value = 4351.3706  # H122: 4354.815 -> 4351.3706
$ws.Cells.Item(122, 9).Value = 4283.6924  # I122: 4290.846 -> 4283.6924
$ws.Cells.Item(122, 11).Value = 12851.0772  # K122: 12872.538 -> 12851.0772
$ws.Cells.Item(122, 13).Value = -10401.0772  # M122: -10422.538 -> -10401.0772

# Row 132
$ws.Cells.Item(132, 8).Value = 5439.1562  # H132: 6260.5386 -> 5439.1562
$ws.Cells.Item(132, 9).Value = 2667.5925  # I132: 2892.6667 -> 2667.5925
$ws.Cells.Item(132, 11).Value = 8002.7775  # K132: 8678.000100000001 -> 8002.7775
$ws.Cells.Item(132, 13).Value = -5472.7775  # M132: -6148.000100000001 -> -5472.7775

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Cells.Item(16, 8).Value = 1996.7675  # H16: 2266.6223 -> 1996.7675
$ws.Cells.Item(16, 9).Value = 1678.1666  # I16: 1806.3158 -> 1678.1666
$ws.Cells.Item(16, 10).Value = 3635.2856  # J16: 4765.4287 -> 3635.2856
$ws.Cells.Item(16, 11).Value = 1678.1666  # K16: 1806.3158 -> 1678.1666
$ws.Cells.Item(16, 12).Value = 3635.2856  # L16: 4765.4287 -> 3635.2856
$ws.Cells.Item(16, 13).Value = -1508.1666  # M16: -1636.3158 -> -1508.1666
$ws.Cells.Item(16, 14).Value = -3975.2856  # N16: -5105.4287 -> -3975.2856

# Row 61
$ws.Cells.Item(61, 8).Value = 2083.423  # H61: 2163.111 -> 2083.423
$ws.Cells.Item(61, 9).Value = 1703.0435  # I61: 1808.5416 -> 1703.0435
$ws.Cells.Item(61, 11).Value = 1703.0435  # K61: 1808.5416 -> 1703.0435
$ws.Cells.Item(61, 13).Value = -1501.0435  # M61: -1606.5416 -> -1501.0435

# Row 100
$ws.Cells.Item(100, 8).Value = 334833.84  # H100: 223722.11 -> 334833.84
$ws.Cells.Item(100, 9).Value = 1667.6666  # I100: 1600 -> 1667.6666
$ws.Cells.Item(100, 10).Value = 668000  # J100: 501374.75 -> 668000
$ws.Cells.Item(100, 11).Value = 1667.6666  # K100: 1600 -> 1667.6666
$ws.Cells.Item(100, 12).Value = 668000  # L100: 501374.75 -> 668000
$ws.Cells.Item(100, 13).Value = -1126.6666  # M100: -1059 -> -1126.6666
$ws.Cells.Item(100, 14).Value = -669082  # N100: -502456.75 -> -669082

# Row 113
$ws.Cells.Item(113, 8).Value = 2083.423  # H113: 2163.111 -> 2083.423
$ws.Cells.Item(113, 9).Value = 1703.0435  # I113: 1808.5416 -> 1703.0435
$ws.Cells.Item(113, 11).Value = 1703.0435  # K113: 1808.5416 -> 1703.0435
$ws.Cells.Item(113, 13).Value = 466.9565  # M113: 361.4584 -> 466.9565

# Row 122
$ws.Cells.Item(122, 8).Value = 4235.385  # H122: 4250.769 -> 4235.385
$ws.Cells.Item(122, 9).Value = 3959.8572  # I122: 4111.385 -> 3959.8572
$ws.Cells.Item(122, 10).Value = 4556.8335  # J122: 4390.154 -> 4556.8335
$ws.Cells.Item(122, 11).Value = 11879.5716  # K122: 12334.155 -> 11879.5716
$ws.Cells.Item(122, 12).Value = 13670.5005  # L122: 13170.462 -> 13670.5005
$ws.Cells.Item(122, 13).Value = -9429.571599999999  # M122: -9884.155000000001 -> -9429.571599999999
$ws.Cells.Item(122, 14).Value = -18570.5005  # N122: -18070.462 -> -18570.5005

# Row 132
$ws.Cells.Item(132, 8).Value = 3315.4324  # H132: 3401.9143 -> 3315.4324
$ws.Cells.Item(132, 9).Value = 3063.3125  # I132: 3147.4 -> 3063.3125
$ws.Cells.Item(132, 11).Value = 9189.9375  # K132: 9442.200000000001 -> 9189.9375
$ws.Cells.Item(132, 13).Value = -6659.9375  # M132: -6912.200000000001 -> -6659.9375

# Row 136
$ws.Cells.Item(136, 8).Value = 2886.5833  # H136: 3139.8064 -> 2886.5833
$ws.Cells.Item(136, 9).Value = 2519.2812  # I136: 2742 -> 2519.2812
$ws.Cells.Item(136, 11).Value = 7557.8436  # K136: 8226 -> 7557.8436
$ws.Cells.Item(136, 13).Value = -5007.8436  # M136: -5676 -> -5007.8436

$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Cells.Item(122, 8).Value = 1864.6666  # H122: 1890 -> 1864.6666
$ws.Cells.Item(122, 9).Value = 1781.5186  # I122: 1794.4814 -> 1781.5186
$ws.Cells.Item(122, 10).Value = 2613  # J122: 2534.75 -> 2613
$ws.Cells.Item(122, 11).Value = 5344.5558  # K122: 5383.4442 -> 5344.5558
$ws.Cells.Item(122, 12).Value = 7839  # L122: 7604.25 -> 7839
$ws.Cells.Item(122, 13).Value = -2894.5558  # M122: -2933.4442 -> -2894.5558
$ws.Cells.Item(122, 14).Value = -12739  # N122: -12504.25 -> -12739

# Row 131
$ws.Cells.Item(131, 8).Value = 146499.5  # H131: 147499.67 -> 146499.5
$ws.Cells.Item(131, 10).Value = 146499.5  # J131: 147499.67 -> 146499.5
$ws.Cells.Item(131, 12).Value = 146499.5  # L131: 147499.67 -> 146499.5
$ws.Cells.Item(131, 14).Value = -156579.5  # N131: -157579.67 -> -156579.5

# Row 136
$ws.Cells.Item(136, 8).Value = 2389.8286  # H136: 2574.3125 -> 2389.8286
$ws.Cells.Item(136, 9).Value = 1245.3636  # I136: 1375.3684 -> 1245.3636
$ws.Cells.Item(136, 11).Value = 3736.0908  # K136: 4126.1052 -> 3736.0908
$ws.Cells.Item(136, 13).Value = -1186.0908  # M136: -1576.1052 -> -1186.0908
